# Weekly price update: a new record is inserted as row 66 (Feria Lagunitas
# de Puerto Montt, Zapallo italiano), pushing all subsequent rows down by
# one. The sheet grows from A1:R138 to A1:R139.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66; Excel shifts rows 66..138 down to 67..139
# and copies formatting (incl. the date number format on column D) from the
# row above, matching the style="2" seen on every D-column cell.
$ws.Rows.Item(66).Insert()

# Populate the new row with the new weekly price record.
$ws.Cells.Item(66, 1).Value  = 4
$ws.Cells.Item(66, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(66, 3).Value  = 'Los Lagos'
$ws.Cells.Item(66, 4).Value  = 44494
$ws.Cells.Item(66, 5).Value  = 10
$ws.Cells.Item(66, 6).Value  = 100112032
$ws.Cells.Item(66, 7).Value  = 'Zapallo italiano'
$ws.Cells.Item(66, 8).Value  = 'Sin especificar'
$ws.Cells.Item(66, 9).Value  = 'Primera'
$ws.Cells.Item(66, 10).Value = 100
$ws.Cells.Item(66, 11).Value = 15000
$ws.Cells.Item(66, 12).Value = 15000
$ws.Cells.Item(66, 13).Value = 15000
$ws.Cells.Item(66, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(66, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(66, 16).Value = 300
$ws.Cells.Item(66, 17).Value = 50
$ws.Cells.Item(66, 18).Value = 'Hortaliza'
